$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.865.30"
$ws.Range("D3").Value = "'1.705.61"
$ws.Range("E3").Value = "'  +0.19%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.30%  "
$ws.Range("D5").Value = "'317.22"
$ws.Range("E5").Value = "'  -0.03%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  -0.37%  "
$ws.Range("D7").Value = "'0.3959"
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("D8").Value = "'0.4063"
$ws.Range("E8").Value = "'  -0.58%  "
$ws.Range("D9").Value = "'1.485"
$ws.Range("E9").Value = "'  -1.44%  "
$ws.Range("D10").Value = "'1.001"
$ws.Range("E10").Value = "'  -0.35%  "
$ws.Range("D11").Value = "'53.51"
$ws.Range("E11").Value = "'  +1.65%  "
$ws.Range("D12").Value = "'0.08812"
$ws.Range("E12").Value = "'  -1.18%  "
$ws.Range("D13").Value = "'26.37"
$ws.Range("E13").Value = "'  +8.17%  "
$ws.Range("D14").Value = "'7.489"
$ws.Range("E14").Value = "'  -2.84%  "
$ws.Range("D15").Value = "'8.123"
$ws.Range("E15").Value = "'  -0.42%  "
$ws.Range("D16").Value = "'0.00001359"
$ws.Range("E16").Value = "'  +1.97%  "
$ws.Range("D17").Value = "'1.689.67"
$ws.Range("E17").Value = "'  -0.96%  "
$ws.Range("D18").Value = "'96.77"
$ws.Range("E18").Value = "'  -3.11%  "
$ws.Range("D19").Value = "'0.07158"
$ws.Range("E19").Value = "'  +0.40%  "
$ws.Range("D20").Value = "'20.98"
$ws.Range("E20").Value = "'  +4.59%  "
$ws.Range("D21").Value = "'7.287"
$ws.Range("E21").Value = "'  +0.82%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "'  -0.68%  "
$ws.Range("D23").Value = "'14.40"
$ws.Range("E23").Value = "'  -1.45%  "
$ws.Range("D24").Value = "'24.864.29"
$ws.Range("E24").Value = "'  +0.29%  "
$ws.Range("D25").Value = "'2.986"
$ws.Range("E25").Value = "'  -3.89%  "
$ws.Range("D26").Value = "'2.334"
$ws.Range("E26").Value = "'  -0.17%  "
$ws.Range("E27").Value = "'  +0.81%  "
$ws.Range("D28").Value = "'6.290"
$ws.Range("E28").Value = "'  +21.11%  "
$ws.Range("D29").Value = "'166.60"
$ws.Range("E29").Value = "'  +0.86%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'145.50"
$ws.Range("E30").Value = "'  +4.15%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'8.394"
$ws.Range("E31").Value = "'  -10.30%  "
$ws.Range("D32").Value = "'2.259"
$ws.Range("E32").Value = "'  +14.87%  "
$ws.Range("D33").Value = "'1.889.55"
$ws.Range("E33").Value = "'  -0.23%  "
$ws.Range("D34").Value = "'0.08784"
$ws.Range("E34").Value = "'  -4.22%  "
$ws.Range("D35").Value = "'0.03193"
$ws.Range("E35").Value = "'  +4.60%  "
$ws.Range("D36").Value = "'7.188"
$ws.Range("E36").Value = "'  -11.34%  "
$ws.Range("D37").Value = "'1.034"
$ws.Range("E37").Value = "'  -4.45%  "
$ws.Range("D38").Value = "'0.2902"
$ws.Range("E38").Value = "'  +3.06%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.8470"
$ws.Range("E39").Value = "'  +8.14%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'10.91"
$ws.Range("E40").Value = "'  -1.62%  "
$ws.Range("D41").Value = "'0.09269"
$ws.Range("E41").Value = "'  -0.47%  "
$ws.Range("D42").Value = "'14.09"
$ws.Range("E42").Value = "'  -3.34%  "
$ws.Range("E43").Value = "'  -0.68%  "
$ws.Range("D44").Value = "'17.50"
$ws.Range("E44").Value = "'  +7.32%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.7456"
$ws.Range("E45").Value = "'  +2.53%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.687"
$ws.Range("E46").Value = "'  +1.87%  "
$ws.Range("D47").Value = "'4.243"
$ws.Range("E47").Value = "'  -0.19%  "
$ws.Range("D48").Value = "'1.391"
$ws.Range("E48").Value = "'  +2.13%  "
$ws.Range("D49").Value = "'0.9999"
$ws.Range("E49").Value = "'  -0.31%  "
$ws.Range("D50").Value = "'141.03"
$ws.Range("E50").Value = "'  -0.14%  "
$ws.Range("D51").Value = "'0.08351"
$ws.Range("E51").Value = "'  +3.53%  "
